# Update "想去人数" (number of people interested) counts on several rows
# across the workbook's sheets, reflecting the latest scrape snapshot.
# Corresponds to commit: "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2866    # 广州·第一届Redamancy动漫游戏嘉年华: 2863 -> 2866
$ws1.Range("F4").Value = 20948   # 广州·2024 CICF×AGF...: 20940 -> 20948
$ws1.Range("F6").Value = 2749    # 广州·南部动漫节: 2746 -> 2749
$ws1.Range("F16").Value = 180    # 广州·2024亚太宠物水族交易会（PSC）国际爬宠展: 179 -> 180

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 125     # 广州·平野宏周粉丝见面会: 124 -> 125
$ws2.Range("F12").Value = 98     # 广州·majiko巡演-2024: 97 -> 98
$ws2.Range("F13").Value = 40     # 广州·「心随歌行」KOKIA 2024 中国巡演: 0 -> 40

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6118    # 广州·「光与夜之恋...」线条大作战主题餐厅: 6117 -> 6118
$ws3.Range("F5").Value = 1553    # 广州·2024《世界之外》x 萌果酱谷子咖啡: 1550 -> 1553

# --- Sheet "全部类型" (All types, aggregated view) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6118    # 广州·「光与夜之恋...」线条大作战主题餐厅: 6117 -> 6118
$ws4.Range("F5").Value = 1553    # 广州·2024《世界之外》x 萌果酱谷子咖啡: 1550 -> 1553
$ws4.Range("F6").Value = 2866    # 广州·第一届Redamancy动漫游戏嘉年华: 2863 -> 2866
$ws4.Range("F8").Value = 20948   # 广州·2024 CICF×AGF...: 20940 -> 20948
$ws4.Range("F12").Value = 125    # 广州·平野宏周粉丝见面会: 124 -> 125
$ws4.Range("F14").Value = 2749   # 广州·南部动漫节: 2746 -> 2749
$ws4.Range("F31").Value = 98     # 广州·majiko巡演-2024: 97 -> 98
$ws4.Range("F32").Value = 180    # 广州·2024亚太宠物水族交易会（PSC）国际爬宠展: 179 -> 180
